# Keywords.xlsx: rename the "Store" keyword group to "Put".
#
# - Worksheet "Store" -> "Put"
# - Keyword name cells "StoreValueIn" -> "PutValueIn"
#   and "StoreTextIn" -> "PutTextIn" (descriptions/arguments text unchanged)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Store")

$ws.Name = "Put"

$ws.Range("A2").Value = "PutValueIn"
$ws.Range("A3").Value = "PutTextIn"
